$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.814.30"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.271.27"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.28"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.59"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.42"
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.33"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.67"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "2.624.01"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.28"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "2.270.55"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.781"
$ws.Range("E18").Value = "  +3.14%  "
$ws.Range("D19").Value = "41.775.33"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.76"
$ws.Range("E20").Value = "  +2.92%  "
$ws.Range("D21").Value = "0.0₃0908"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.93"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "243.36"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.93"
$ws.Range("E26").Value = "  +2.33%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.98"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.54"
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("E30").Value = "  -5.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.93"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.12"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.24"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.01"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.86"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.79"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.93"
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("D43").Value = "2.018.15"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.59"
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.38"
$ws.Range("E46").Value = "  +2.97%  "
$ws.Range("E47").Value = "  +7.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.89"
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.17"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.08"
$ws.Range("E50").Value = "  +3.19%  "
$ws.Range("E51").Value = "  -1.48%  "
